# Fruta / hortaliza, semanal
# Insert one new weekly data point at row 18 (pushing the existing
# rows 18-68 down to 19-69) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 18, shifting rows 18:68 down to 19:69.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value = 44707
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100102
$ws.Range("H18").Value = "Cítricos"
$ws.Range("I18").Value = 100102006
$ws.Range("J18").Value = "Pomelo"
$ws.Range("K18").Value = "Start Ruby"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 280
$ws.Range("N18").Value = 7500
$ws.Range("O18").Value = 7500
$ws.Range("P18").Value = 7500
$ws.Range("Q18").Value = "$/caja 14 kilos"
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 536
$ws.Range("T18").Value = 14
